$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text representation
# (values like "135.44" or "0.0970" must stay as text, matching the source
# workbook where these cells are inline strings, not numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '60.182.34'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").Value = '2.407.46'
$ws.Range("E3").Value = '  -0.67%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '559.64'
$ws.Range("E5").Value = '  +1.29%  '

$ws.Range("D6").Value = '135.44'
$ws.Range("E6").Value = '  -1.92%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  +0.32%  '

$ws.Range("D9").Value = '0.106'
$ws.Range("E9").Value = '  -0.57%  '

$ws.Range("D10").Value = '5.61'
$ws.Range("E10").Value = '  -1.71%  '

$ws.Range("E11").Value = '  +0.38%  '

$ws.Range("E12").Value = '  -1.89%  '

$ws.Range("D13").Value = '24.74'
$ws.Range("E13").Value = '  -2.11%  '

$ws.Range("D14").Value = '2.839.47'
$ws.Range("E14").Value = '  -0.57%  '

$ws.Range("D15").Value = '60.103.49'
$ws.Range("E15").Value = '  +0.40%  '

$ws.Range("E16").Value = '  -0.10%  '

$ws.Range("D17").Value = '2.341.41'
$ws.Range("E17").Value = '  -4.00%  '

$ws.Range("D18").Value = '11.23'
$ws.Range("E18").Value = '  -0.89%  '

$ws.Range("D19").Value = '4.53'
$ws.Range("E19").Value = '  +2.92%  '

$ws.Range("D20").Value = '326.20'
$ws.Range("E20").Value = '  -1.52%  '

$ws.Range("E21").Value = '  +1.62%  '

$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").Value = '64.71'
$ws.Range("E23").Value = '  -2.35%  '

$ws.Range("D24").Value = '0.171'
$ws.Range("E24").Value = '  -0.41%  '

$ws.Range("E25").Value = '  -2.37%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("E27").Value = '  +0.78%  '

$ws.Range("E28").Value = '  +1.46%  '

$ws.Range("D29").Value = '0.0₃0769'
$ws.Range("E29").Value = '  -1.23%  '

$ws.Range("D30").Value = '170.78'
$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("D31").Value = '6.12'
$ws.Range("E31").Value = '  -0.34%  '

$ws.Range("E32").Value = '  +5.31%  '

$ws.Range("E33").Value = '  -2.32%  '

$ws.Range("D34").Value = '18.41'
$ws.Range("E34").Value = '  -1.51%  '

$ws.Range("D35").Value = '1.34'
$ws.Range("E35").Value = '  +3.40%  '

$ws.Range("D38").Value = '4.18'
$ws.Range("E38").Value = '  -1.02%  '

$ws.Range("D39").Value = '325.07'
$ws.Range("E39").Value = '  +3.46%  '

$ws.Range("E40").Value = '  -0.81%  '

$ws.Range("E41").Value = '  -2.76%  '

$ws.Range("D42").Value = '148.76'
$ws.Range("E42").Value = '  +6.90%  '

$ws.Range("D43").Value = '3.57'
$ws.Range("E43").Value = '  -3.45%  '

$ws.Range("D44").Value = '0.0970'
$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("D45").Value = '19.97'

$ws.Range("E46").Value = '  -0.86%  '

$ws.Range("E47").Value = '  -0.43%  '

$ws.Range("E48").Value = '  -1.70%  '

$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("E50").Value = '  -1.44%  '

$ws.Range("D51").Value = '4.66'
$ws.Range("E51").Value = '  -0.85%  '
